# Add a new "Spiral Matrix" coding-question row (row 57) to the hint sheet,
# copying the formatting of the row above it (row 56) and filling in the
# new question's data, then move the active selection to the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles) of the existing last data row (56) down into
# the new row (57) so number formats / alignment / fonts match the rest of
# the table.
$ws.Range("A56:G56").Copy()
$ws.Range("A57:G57").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the values for the new row.
$ws.Range("A57").Value2 = 54
$ws.Range("B57").Value2 = 45666
$ws.Range("C57").Value2 = "Spiral Matrix. Leetcode"
$ws.Range("F57").Value2 = $ws.Range("F56").Value2
$ws.Range("G57").Value2 = $ws.Range("G56").Value2

# Bold the trailing "Leetcode" portion of the problem-statement cell, same
# styling convention used by the other rows in the sheet.
$chars = $ws.Range("C57").Characters(16, 8)
$chars.Font.Bold = $true
$chars.Font.Name = "Calibri"
$chars.Font.Size = 11

# Update the current selection to the newly added row, matching where the
# author's cursor ended up after the edit.
$ws.Range("E57").Select()
